$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new metric data point as row 42 (directly below the last existing row, 41)
$newRow = 42

$ws.Cells.Item($newRow, 1).Value = "2025-04-29 04:13:51"
$ws.Cells.Item($newRow, 2).Value = 112
